$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Heading3 title: "Freelance Consultant / Technology Advisor"
#    -> "Chief Technology Officer at ARTPIX jewelry visualizations"
#    (split into runs like the sibling "Senior Software Engineer at ..."
#    headings use, and drop the now-redundant pPr/rPr color override)
# -----------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Freelance Consultant / Technology Advisor") {
        $titlePara = $p
        break
    }
}

$titleParaRange = $titlePara.Range

$titleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">Chief Technology Officer </w:t></w:r><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="auto"/></w:rPr><w:t>at</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr/><w:t>ARTPIX jewelry visualizations</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$titleParaRange.InsertXML($titleXml)

# -----------------------------------------------------------------------
# 2) Date line directly below it: "2021 - 2023" -> "2023 - 2025"
#    (split into 4 runs the same way the diff shows, keeping the
#    paragraph's existing italic/grey formatting). Anchor on the
#    paragraph immediately following the title we just edited, since
#    the same "2021 - 2023" text also appears later in the CV.
# -----------------------------------------------------------------------
$datePara = $titlePara.Next()

$dateParaRange = $datePara.Range

$dateRpr = '<w:rPr><w:i w:val="1"/><w:iCs w:val="1"/><w:color w:val="808080" w:themeColor="background1" w:themeTint="FF" w:themeShade="80"/></w:rPr>'

$dateXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/>$dateRpr</w:pPr><w:r>$dateRpr<w:t>202</w:t></w:r><w:r>$dateRpr<w:t>3</w:t></w:r><w:r>$dateRpr<w:t xml:space="preserve"> - 202</w:t></w:r><w:r>$dateRpr<w:t>5</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

[void]$dateParaRange.InsertXML($dateXml)

Write-Host "Title now: [$($titlePara.Range.Text)]"
Write-Host "Date now: [$($datePara.Range.Text)]"
